$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2535.6428
$ws.Range("I132").Value = 1780
$ws.Range("J132").Value = 4424.75
$ws.Range("K132").Value = 5340
$ws.Range("L132").Value = 13274.25
$ws.Range("M132").Value = -2810
$ws.Range("N132").Value = -18334.25
$ws.Range("H138").Value = 3641.1855
$ws.Range("I138").Value = 1650.5834
$ws.Range("J138").Value = 3922.2117
$ws.Range("K138").Value = 4951.7502
$ws.Range("L138").Value = 11766.6351
$ws.Range("M138").Value = 188.2497999999996
$ws.Range("N138").Value = -22046.6351
$ws.Range("H141").Value = 776.7959
$ws.Range("I141").Value = 776.7959
$ws.Range("K141").Value = 2330.3877
$ws.Range("M141").Value = 2849.6123

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6228.0566
$ws.Range("I32").Value = 4152.5317
$ws.Range("K32").Value = 4152.5317
$ws.Range("M32").Value = -3865.5317
$ws.Range("H44").Value = 27374.5
$ws.Range("J44").Value = 29999.334
$ws.Range("L44").Value = 29999.334
$ws.Range("N44").Value = -30975.334
$ws.Range("H45").Value = 211561.9
$ws.Range("I45").Value = 221523.8
$ws.Range("J45").Value = 201600
$ws.Range("K45").Value = 221523.8
$ws.Range("L45").Value = 201600
$ws.Range("M45").Value = -221146.8
$ws.Range("N45").Value = -202354
$ws.Range("H55").Value = 11999.5
$ws.Range("J55").Value = 19999
$ws.Range("L55").Value = 19999
$ws.Range("N55").Value = -20629
$ws.Range("H61").Value = 2877.5366
$ws.Range("I61").Value = 2082.5173
$ws.Range("K61").Value = 2082.5173
$ws.Range("M61").Value = -1870.5173
$ws.Range("H97").Value = 2583.2188
$ws.Range("I97").Value = 1885.5862
$ws.Range("K97").Value = 1885.5862
$ws.Range("M97").Value = -1389.5862
$ws.Range("H136").Value = 2877.5366
$ws.Range("I136").Value = 2082.5173
$ws.Range("K136").Value = 6247.5519
$ws.Range("M136").Value = -3697.5519

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 2302.6667
$ws.Range("J5").Value = 1000
$ws.Range("L5").Value = 1000
$ws.Range("N5").Value = -1226
$ws.Range("H105").Value = 2082.818
$ws.Range("I105").Value = 1841.2858
$ws.Range("K105").Value = 1841.2858
$ws.Range("M105").Value = -94.28580000000011
$ws.Range("H134").Value = 1782.7115
$ws.Range("I134").Value = 1782.7115
$ws.Range("K134").Value = 5348.1345
$ws.Range("M134").Value = -2813.1345

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15292.975
$ws.Range("I31").Value = 20612.426
$ws.Range("J31").Value = 4244.885
$ws.Range("K31").Value = 20612.426
$ws.Range("L31").Value = 4244.885
$ws.Range("M31").Value = -20317.426
$ws.Range("N31").Value = -4834.885
$ws.Range("H34").Value = 15292.975
$ws.Range("I34").Value = 20612.426
$ws.Range("J34").Value = 4244.885
$ws.Range("K34").Value = 20612.426
$ws.Range("L34").Value = 4244.885
$ws.Range("M34").Value = -20410.426
$ws.Range("N34").Value = -4648.885
$ws.Range("H57").Value = 46666.332
$ws.Range("I57").Value = 54999
$ws.Range("J57").Value = 42500
$ws.Range("K57").Value = 54999
$ws.Range("L57").Value = 42500
$ws.Range("M57").Value = -54439
$ws.Range("N57").Value = -43620
$ws.Range("H58").Value = 3043.7954
$ws.Range("I58").Value = 2770.2307
$ws.Range("K58").Value = 2770.2307
$ws.Range("M58").Value = -2567.2307
$ws.Range("H99").Value = 4069.0557
$ws.Range("I99").Value = 3954.7273
$ws.Range("K99").Value = 3954.7273
$ws.Range("M99").Value = -2456.7273
$ws.Range("H126").Value = 4069.0557
$ws.Range("I126").Value = 3954.7273
$ws.Range("K126").Value = 11864.1819
$ws.Range("M126").Value = -9394.1819
$ws.Range("H132").Value = 2262.2292
$ws.Range("I132").Value = 1789.875
$ws.Range("K132").Value = 5369.625
$ws.Range("M132").Value = -2839.625
$ws.Range("H134").Value = 16312.31
$ws.Range("I134").Value = 5206.057
$ws.Range("J134").Value = 71843.57000000001
$ws.Range("K134").Value = 15618.171
$ws.Range("L134").Value = 215530.71
$ws.Range("M134").Value = -13083.171
$ws.Range("N134").Value = -220600.71
$ws.Range("H136").Value = 3043.7954
$ws.Range("I136").Value = 2770.2307
$ws.Range("K136").Value = 8310.6921
$ws.Range("M136").Value = -5760.6921

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5276742.5
$ws.Range("I4").Value = 5276742.5
$ws.Range("K4").Value = 15830227.5
$ws.Range("M4").Value = -15830115.5
$ws.Range("H56").Value = 9165
$ws.Range("I56").Value = 9165
$ws.Range("K56").Value = 9165
$ws.Range("M56").Value = -8635
$ws.Range("H107").Value = 1576.1818
$ws.Range("I107").Value = 836
$ws.Range("K107").Value = 2508
$ws.Range("M107").Value = -588
$ws.Range("H132").Value = 1101.8049
$ws.Range("I132").Value = 1058.7567
$ws.Range("K132").Value = 9528.810299999999
$ws.Range("M132").Value = -6998.810299999999
$ws.Range("H139").Value = 2221
$ws.Range("I139").Value = 2123.625
$ws.Range("J139").Value = 3000
$ws.Range("K139").Value = 6370.875
$ws.Range("L139").Value = 9000
$ws.Range("M139").Value = -1230.875
$ws.Range("N139").Value = -19280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 4699.4
$ws.Range("I24").Value = 4499
$ws.Range("J24").Value = 5000
$ws.Range("K24").Value = 4499
$ws.Range("L24").Value = 5000
$ws.Range("M24").Value = -4326
$ws.Range("N24").Value = -5346
$ws.Range("H33").Value = 10001
$ws.Range("J33").Value = 10001
$ws.Range("L33").Value = 10001
$ws.Range("N33").Value = -10505
$ws.Range("H36").Value = 10839.333
$ws.Range("J36").Value = 15001
$ws.Range("L36").Value = 15001
$ws.Range("N36").Value = -15971
$ws.Range("H52").Value = 29500
$ws.Range("J52").Value = 29000
$ws.Range("L52").Value = 29000
$ws.Range("N52").Value = -29518
$ws.Range("H58").Value = 54078.5
$ws.Range("J58").Value = 54078.5
$ws.Range("L58").Value = 54078.5
$ws.Range("N58").Value = -54632.5
$ws.Range("H132").Value = 3455.4285
$ws.Range("I132").Value = 3480.7932
$ws.Range("J132").Value = 3332.8333
$ws.Range("K132").Value = 10442.3796
$ws.Range("L132").Value = 9998.499899999999
$ws.Range("M132").Value = -7912.3796
$ws.Range("N132").Value = -15058.4999
$ws.Range("H136").Value = 37133.547
$ws.Range("J136").Value = 37133.547
$ws.Range("L136").Value = 111400.641
$ws.Range("N136").Value = -116500.641

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2707.842
$ws.Range("I16").Value = 1825.4
$ws.Range("J16").Value = 6017
$ws.Range("K16").Value = 1825.4
$ws.Range("L16").Value = 6017
$ws.Range("M16").Value = -1655.4
$ws.Range("N16").Value = -6357
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H61").Value = 39430.586
$ws.Range("I61").Value = 60207.668
$ws.Range("J61").Value = 5431.727
$ws.Range("K61").Value = 60207.668
$ws.Range("L61").Value = 5431.727
$ws.Range("M61").Value = -60005.668
$ws.Range("N61").Value = -5835.727
$ws.Range("H100").Value = 49998.5
$ws.Range("I100").Value = 49998
$ws.Range("K100").Value = 49998
$ws.Range("M100").Value = -49457
$ws.Range("H106").Value = 19999.5
$ws.Range("J106").Value = 19999.5
$ws.Range("L106").Value = 19999.5
$ws.Range("N106").Value = -22523.5
$ws.Range("H113").Value = 39430.586
$ws.Range("I113").Value = 60207.668
$ws.Range("J113").Value = 5431.727
$ws.Range("K113").Value = 60207.668
$ws.Range("L113").Value = 5431.727
$ws.Range("M113").Value = -58037.668
$ws.Range("N113").Value = -9771.726999999999
$ws.Range("H132").Value = 4318.5713
$ws.Range("I132").Value = 3747.2354
$ws.Range("J132").Value = 6746.75
$ws.Range("K132").Value = 11241.7062
$ws.Range("L132").Value = 20240.25
$ws.Range("M132").Value = -8711.706200000001
$ws.Range("N132").Value = -25300.25
$ws.Range("H136").Value = 3153.5806
$ws.Range("I136").Value = 2866.4644
$ws.Range("K136").Value = 8599.393199999999
$ws.Range("M136").Value = -6049.393199999999
$ws.Range("H140").Value = 187239.25
$ws.Range("J140").Value = 187239.25
$ws.Range("L140").Value = 187239.25
$ws.Range("N140").Value = -197599.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 40307.5
$ws.Range("J105").Value = 40307.5
$ws.Range("L105").Value = 40307.5
$ws.Range("N105").Value = -47295.5
$ws.Range("H122").Value = 1677.1666
$ws.Range("I122").Value = 1549.8055
$ws.Range("J122").Value = 2441.3333
$ws.Range("K122").Value = 4649.416499999999
$ws.Range("L122").Value = 7323.999899999999
$ws.Range("M122").Value = -2199.416499999999
$ws.Range("N122").Value = -12223.9999
$ws.Range("H126").Value = 2327.6316
$ws.Range("I126").Value = 1857.8125
$ws.Range("K126").Value = 5573.4375
$ws.Range("M126").Value = -3103.4375
$ws.Range("H136").Value = 1550.4861
$ws.Range("I136").Value = 1530.3019
$ws.Range("K136").Value = 4590.905699999999
$ws.Range("M136").Value = -2040.905699999999
